# Add 2022-Q4 data
#
# 1. "总计" (summary) sheet: insert a new row 2 for the 2022-Q4 quarter,
#    pushing the existing quarters down by one row, and bump the running
#    index in column A for every row that moved.
# 2. Insert a brand-new "2022-Q4" worksheet right after "总计", containing
#    the per-fund holdings detail for that quarter.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# --- Step 1: update the "总计" sheet -------------------------------------

$total.Rows.Item(2).Insert()

# Insert() copies the row-above formatting into the new blank row; put the
# data cells back to the unformatted style used by every other data row,
# then restore column A's header-like style (bold/bordered, like the rows
# below it) by copying formats from the row underneath.
$total.Range("B2:D2").Style = "Normal"
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 0.32

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7

# --- Step 2: insert the new "2022-Q4" sheet right after "总计" -----------

$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Match the bold/bordered header style used on every other quarter sheet.
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Column A (the running index) uses the same centered/bordered style as
# column A on every other sheet in the workbook.
$total.Range("A2").Copy()
$q4.Range("A2:A11").PasteSpecial(-4122)

$rows = @(
    @(0, "550001", "信诚四季红混合",                         "4.65", "82.79", "2.74", "0.1274", 8),
    @(1, "001656", "农银汇理中国优势灵活配置混合",             "2.01", "91.26", "2.99", "0.0601", 2),
    @(2, "001060", "前海开源高端装备制造灵活配置混合",         "0.90", "88.69", "6.20", "0.0558", 1),
    @(3, "011284", "中信保诚龙腾精选混合",                     "1.08", "83.70", "2.73", "0.0295", 8),
    @(4, "004750", "广发鑫和灵活配置混合A",                    "2.29", "21.57", "1.19", "0.0273", 2),
    @(5, "004751", "广发鑫和灵活配置混合C",                    "0.68", "21.57", "1.19", "0.0081", 2),
    @(6, "006123", "中融高股息精选混合A",                      "0.22", "87.21", "3.15", "0.0069", 7),
    @(7, "006124", "中融高股息精选混合C",                      "0.16", "87.21", "3.15", "0.0050", 7),
    @(8, "001412", "德邦鑫星价值灵活配置混合A",                "0.13", "35.79", "1.86", "0.0024", 9),
    @(9, "002112", "德邦鑫星价值灵活配置混合C",                "0.02", "35.79", "1.86", "0.0004", 9)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    # Fund code / size / position figures are stored as text in the source
    # data (leading zeros, fixed decimal places), so force text with a
    # leading quote instead of letting them coerce to numbers, then drop
    # the quote-prefix formatting it implies so the cell is plain text
    # with no explicit style, matching the source rows.
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $q4.Range("B" + $r + ":G" + $r).Style = "Normal"
    $r = $r + 1
}
